$d = $word.ActiveDocument
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 1. Title (Heading1) and the later bold repeat of the title
Replace-Text "Play Fire Opals Slot for Free - Review and Ratings" "Play Fire Opals for Free - A Mystical Slot Game"

# 2. "What we like" bullets - reshuffled content, replace using unique anchors
Replace-Text "720 paylines for increased volatility" "Unique prism-like structure with 720 paylines"
Replace-Text "Winning combinations from both left to right and right to left" "Special symbols and free spin mode"
Replace-Text "Special symbols and free spin mode with up to 260 spins" "Visually appealing graphics and symbols"
Replace-Text "Unique prism-like structure and dreamy Hawaiian setting" "Winning combinations from left to right and right to left"

# 3. "What we don't like" bullets
Replace-Text "Tribal soundtrack might not fit the dreamy atmosphere" "Playful sound design may not match the game's atmosphere"
Replace-Text "Limited number of special symbols" "Not the only Hawaiian-themed slot available"

# 5. Meta description (italic)
Replace-Text "Explore Fire Opals online slot and play for free with bonus rounds, wilds, and free spins. Read a full review and ratings, compare to other Hawaiian-themed slots." "Read our review of Fire Opals, a mystical slot game set in Hawaii. Play it for free and enjoy its unique features."
